# Update example template with VERSION sheet
#
# 1. Add three vendor names ("Starbucks", "Amazon", "AWS") to the Setup
#    sheet's "Vendors" list (column F, rows 2-4).
# 2. Append a new "VERSION" worksheet at the end of the workbook with a
#    Property/Value table describing the generated build.

$wb = $excel.ActiveWorkbook

# --- 1. Setup sheet: populate the Vendors column -------------------------
$setup = $wb.Worksheets.Item("Setup")
$setup.Range("F2").Value = "Starbucks"
$setup.Range("F3").Value = "Amazon"
$setup.Range("F4").Value = "AWS"

# --- 2. Add the new VERSION sheet as the last tab -------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$versionSheet = $wb.Worksheets.Add($null, $lastSheet)
$versionSheet.Name = "VERSION"

$versionSheet.Range("A1").Value = "Property"
$versionSheet.Range("B1").Value = "Value"

$versionSheet.Range("A2").Value = "Version ID"
$versionSheet.Range("B2").Value = "1.0.0-java"

$versionSheet.Range("A3").Value = "Git SHA"
$versionSheet.Range("B3").Value = "fa8999fc0a5c37c5e114a44cb01475b5da278394"

$versionSheet.Range("A4").Value = "Generated At"
$versionSheet.Range("B4").Value = "12/22/2025, 10:07:03 AM"

# Column widths matching the authored sheet (A=20, B=50 "characters" in the
# raw OOXML <col width="..."/> units). Excel's ColumnWidth COM property is
# expressed in a slightly different unit that adds ~5/6 of a character as
# padding, so back that padding out to land on the exact authored width.
$versionSheet.Columns.Item(1).ColumnWidth = 20 - 5/6
$versionSheet.Columns.Item(2).ColumnWidth = 50 - 5/6
